# Auto-update draw results: append the 2025-12-09 "Pick 4" draw as a new
# row at the bottom of the Results sheet (row 84), matching the layout of
# every prior row: Date | Game | Phase | Result | InsertedAt.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 84

# Columns A ("2025-12-09") and C ("251209") look like a date / a plain
# number to Excel's input parser, so a bare assignment would silently be
# converted into a date serial / numeric value instead of staying text
# (exactly like every other row in this sheet, which stores these columns
# as text). Mark those two cells as Text ("@") before writing into them so
# the value sticks as a literal string, matching the rest of the column.
$ws.Cells.Item($newRow, 1).NumberFormat = "@"
$ws.Cells.Item($newRow, 3).NumberFormat = "@"

$ws.Cells.Item($newRow, 1).Value2 = "2025-12-09"
$ws.Cells.Item($newRow, 2).Value2 = "Pick 4"
$ws.Cells.Item($newRow, 3).Value2 = "251209"
$ws.Cells.Item($newRow, 4).Value2 = "3-6-8-0"
$ws.Cells.Item($newRow, 5).Value2 = "2025-12-09T21:42:45.616+04:00"

# The sheet keeps a "numbers stored as text" ignored-error marker over the
# whole data range; extend it to cover the newly appended row too.
$dataRange = $ws.Range("A1:E84")
$errs = $dataRange.Errors
$errs.Item(3).Ignore = $true
